$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.971.96'
$ws.Range("E2").Value = '  -2.74%  '
$ws.Range("D3").Value = '1.794.01'
$ws.Range("E3").Value = '  -3.17%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = '307.85'
$ws.Range("E5").Value = '  -2.34%  '
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  +0.34%  '
$ws.Range("D7").Value = '0.4179'
$ws.Range("E7").Value = '  -2.99%  '
$ws.Range("D8").Value = '0.3555'
$ws.Range("E8").Value = '  -4.00%  '
$ws.Range("D9").Value = '0.07079'
$ws.Range("E9").Value = '  -3.63%  '
$ws.Range("D10").Value = '0.8420'
$ws.Range("E10").Value = '  -3.76%  '
$ws.Range("D11").Value = '20.13'
$ws.Range("E11").Value = '  -4.16%  '
$ws.Range("D12").Value = '1.841.05'
$ws.Range("E12").Value = '  -2.30%  '
$ws.Range("D13").Value = '5.287'
$ws.Range("E13").Value = '  -2.94%  '
$ws.Range("D14").Value = '6.342'
$ws.Range("E14").Value = '  -3.74%  '
$ws.Range("D15").Value = '0.06743'
$ws.Range("E15").Value = '  -2.93%  '
$ws.Range("D16").Value = '1.008'
$ws.Range("E16").Value = '  +0.51%  '
$ws.Range("D17").Value = '79.51'
$ws.Range("E17").Value = '  -2.19%  '
$ws.Range("D18").Value = '0.000008696'
$ws.Range("E18").Value = '  -4.25%  '
$ws.Range("D19").Value = '1.004'
$ws.Range("E19").Value = '  +0.40%  '
$ws.Range("D20").Value = '15.04'
$ws.Range("E20").Value = '  -3.24%  '
$ws.Range("D21").Value = '26.976.46'
$ws.Range("E21").Value = '  -3.05%  '
$ws.Range("D22").Value = '5.046'
$ws.Range("E22").Value = '  -0.82%  '
$ws.Range("D23").Value = '11.05'
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("D24").Value = '2.022.99'
$ws.Range("E24").Value = '  -3.49%  '
$ws.Range("D25").Value = '1.938'
$ws.Range("E25").Value = '  -1.51%  '
$ws.Range("D26").Value = '152.70'
$ws.Range("E26").Value = '  -1.48%  '
$ws.Range("D27").Value = '18.09'
$ws.Range("E27").Value = '  -2.65%  '
$ws.Range("D28").Value = '4.987'
$ws.Range("E28").Value = '  -6.17%  '
$ws.Range("D29").Value = '112.76'
$ws.Range("E30").Value = '  -12.00%  '
$ws.Range("D31").Value = '0.08901'
$ws.Range("E31").Value = '  -0.48%  '
$ws.Range("D32").Value = '0.7133'
$ws.Range("E32").Value = '  -8.80%  '
$ws.Range("D33").Value = '2.843'
$ws.Range("E33").Value = '  -4.44%  '
$ws.Range("D34").Value = '4.288'
$ws.Range("E34").Value = '  -6.82%  '
$ws.Range("D35").Value = '1.004'
$ws.Range("E35").Value = '  +0.36%  '
$ws.Range("D36").Value = '1.066'
$ws.Range("E36").Value = '  -8.19%  '
$ws.Range("D37").Value = '1.072'
$ws.Range("E37").Value = '  -3.46%  '
$ws.Range("D38").Value = '0.01897'
$ws.Range("E38").Value = '  -3.14%  '
$ws.Range("D39").Value = '0.05097'
$ws.Range("E39").Value = '  -5.98%  '
$ws.Range("D40").Value = '0.1620'
$ws.Range("E40").Value = '  -3.58%  '
$ws.Range("D41").Value = '0.4927'
$ws.Range("E41").Value = '  -5.47%  '
$ws.Range("D42").Value = '2.574'
$ws.Range("E42").Value = '  -9.30%  '
$ws.Range("D43").Value = '5.993'
$ws.Range("E43").Value = '  -11.12%  '
$ws.Range("D44").Value = '8.009'
$ws.Range("E44").Value = '  -7.29%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '10.21'
$ws.Range("E45").Value = '  -4.13%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '1.004'
$ws.Range("E46").Value = '  +0.35%  '
$ws.Range("D47").Value = '104.15'
$ws.Range("D48").Value = '0.06301'
$ws.Range("E48").Value = '  -4.03%  '
$ws.Range("D49").Value = '0.4502'
$ws.Range("E49").Value = '  -5.90%  '
$ws.Range("E50").Value = '  -4.98%  '
$ws.Range("D51").Value = '61.87'
$ws.Range("E51").Value = '  -4.72%  '
